$d = $word.ActiveDocument

$d.Content.Find.Execute("26×98=2548", $true, $false, $false, $false, $false, $true, 1, $false, "15×11=165", 2) | Out-Null
$d.Content.Find.Execute("39×68=2652", $true, $false, $false, $false, $false, $true, 1, $false, "96×37=3552", 2) | Out-Null
$d.Content.Find.Execute("58×79=4582", $true, $false, $false, $false, $false, $true, 1, $false, "56×98=5488", 2) | Out-Null
$d.Content.Find.Execute("73×81=5913", $true, $false, $false, $false, $false, $true, 1, $false, "35×87=3045", 2) | Out-Null
$d.Content.Find.Execute("89×32=2848", $true, $false, $false, $false, $false, $true, 1, $false, "34×24=816", 2) | Out-Null
$d.Content.Find.Execute("29×25=725", $true, $false, $false, $false, $false, $true, 1, $false, "74×44=3256", 2) | Out-Null
$d.Content.Find.Execute("26×33=858", $true, $false, $false, $false, $false, $true, 1, $false, "44×70=3080", 2) | Out-Null
$d.Content.Find.Execute("68×36=2448", $true, $false, $false, $false, $false, $true, 1, $false, "82×90=7380", 2) | Out-Null
$d.Content.Find.Execute("51×61=3111", $true, $false, $false, $false, $false, $true, 1, $false, "40×81=3240", 2) | Out-Null
$d.Content.Find.Execute("55×98=5390", $true, $false, $false, $false, $false, $true, 1, $false, "11×43=473", 2) | Out-Null
$d.Content.Find.Execute("64×57=3648", $true, $false, $false, $false, $false, $true, 1, $false, "72×71=5112", 2) | Out-Null
$d.Content.Find.Execute("91×84=7644", $true, $false, $false, $false, $false, $true, 1, $false, "79×28=2212", 2) | Out-Null
$d.Content.Find.Execute("50×93=4650", $true, $false, $false, $false, $false, $true, 1, $false, "75×66=4950", 2) | Out-Null
$d.Content.Find.Execute("26×24=624", $true, $false, $false, $false, $false, $true, 1, $false, "41×49=2009", 2) | Out-Null
$d.Content.Find.Execute("12×90=1080", $true, $false, $false, $false, $false, $true, 1, $false, "90×75=6750", 2) | Out-Null
$d.Content.Find.Execute("82×78=6396", $true, $false, $false, $false, $false, $true, 1, $false, "49×67=3283", 2) | Out-Null
$d.Content.Find.Execute("64×25=1600", $true, $false, $false, $false, $false, $true, 1, $false, "64×96=6144", 2) | Out-Null
$d.Content.Find.Execute("34×67=2278", $true, $false, $false, $false, $false, $true, 1, $false, "88×60=5280", 2) | Out-Null
$d.Content.Find.Execute("13×37=481", $true, $false, $false, $false, $false, $true, 1, $false, "60×30=1800", 2) | Out-Null
$d.Content.Find.Execute("36×66=2376", $true, $false, $false, $false, $false, $true, 1, $false, "79×61=4819", 2) | Out-Null
$d.Content.Find.Execute("25×52=1300", $true, $false, $false, $false, $false, $true, 1, $false, "60×74=4440", 2) | Out-Null
$d.Content.Find.Execute("23×43=989", $true, $false, $false, $false, $false, $true, 1, $false, "52×51=2652", 2) | Out-Null
$d.Content.Find.Execute("33×30=990", $true, $false, $false, $false, $false, $true, 1, $false, "81×67=5427", 2) | Out-Null
$d.Content.Find.Execute("31×18=558", $true, $false, $false, $false, $false, $true, 1, $false, "86×27=2322", 2) | Out-Null
$d.Content.Find.Execute("76×34=2584", $true, $false, $false, $false, $false, $true, 1, $false, "12×81=972", 2) | Out-Null
